# feat: add 2022-Q1 data
#
# Insert a new "2022-Q1" sheet between "2021-Q4" and "总计" (fund holdings
# for the quarter), and rebuild "总计" with a new summary row for 2022-Q1
# prepended above the existing 2021-Q4 summary row.

$wb = $excel.ActiveWorkbook

# Writes a value as literal text (shared-string), even when it looks like
# a number (e.g. "14.71", "010583"), without leaving any NumberFormat /
# style residue behind on the cell or in the workbook's style table.
function Set-TextValue($range, [string]$text) {
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163) # xlPasteValues
}

$q4 = $wb.Worksheets.Item("2021-Q4")

# Detach the existing "总计" sheet now so that re-adding it after the new
# "2022-Q1" sheet gives the trio the sheetId / tab order from the diff:
# 2021-Q4 (1), 2022-Q1 (2), 总计 (3).
$old_zj = $wb.Worksheets.Item("总计")
$old_zj.Delete()

# ----------------------------------------------------------------------
# New sheet: 2022-Q1
# ----------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

Set-TextValue $q1.Range("B1") "基金代码"
Set-TextValue $q1.Range("C1") "基金名称"
Set-TextValue $q1.Range("D1") "基金规模"
Set-TextValue $q1.Range("E1") "股票总仓位"
Set-TextValue $q1.Range("F1") "仓位占比"
Set-TextValue $q1.Range("G1") "持有市值(亿元)"
Set-TextValue $q1.Range("H1") "仓位排名"

$rows = @(
    @("010583", "富国蓝筹精选股票（QDII）美元",   "14.71", "94.57", "5.06", "0.7443", 3),
    @("007455", "富国蓝筹精选股票（QDII）人民币", "14.71", "94.57", "5.06", "0.7443", 3),
    @("457001", "国富亚洲机会股票 (QDII)",         "5.93",  "77.36", "2.76", "0.1637", 8),
    @("100055", "富国全球科技互联网股票(QDII)",    "3.01",  "70.87", "4.02", "0.1210", 6)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $row = $i + 2
    $data = $rows[$i]

    $q1.Range("A$row").Value = $i

    Set-TextValue $q1.Range("B$row") $data[0]
    Set-TextValue $q1.Range("C$row") $data[1]
    Set-TextValue $q1.Range("D$row") $data[2]
    Set-TextValue $q1.Range("E$row") $data[3]
    Set-TextValue $q1.Range("F$row") $data[4]
    Set-TextValue $q1.Range("G$row") $data[5]

    $q1.Range("H$row").Value = $data[6]
}

# Match the bold/bordered header style used on "2021-Q4" (style index 2)
# by copying formats from that sheet's header / index cells.
$q4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122) # xlPasteFormats

$q4.Range("A2").Copy()
$q1.Range("A2:A5").PasteSpecial(-4122) # xlPasteFormats

# ----------------------------------------------------------------------
# Rebuilt sheet: 总计
# ----------------------------------------------------------------------
$zj = $wb.Worksheets.Add($null, $q1)
$zj.Name = "总计"

Set-TextValue $zj.Range("B1") "日期"
Set-TextValue $zj.Range("C1") "持有数量(只)"
Set-TextValue $zj.Range("D1") "持有市值(亿元)"

$zj.Range("A2").Value = 0
Set-TextValue $zj.Range("B2") "2022-Q1"
$zj.Range("C2").Value = 4
$zj.Range("D2").Value = 1.77

$zj.Range("A3").Value = 1
Set-TextValue $zj.Range("B3") "2021-Q4"
$zj.Range("C3").Value = 3
$zj.Range("D3").Value = 1.8

$q4.Range("B1:D1").Copy()
$zj.Range("B1:D1").PasteSpecial(-4122) # xlPasteFormats

$q4.Range("A2").Copy()
$zj.Range("A2:A3").PasteSpecial(-4122) # xlPasteFormats

# Restore the original active sheet ("2021-Q4"), since creating new
# sheets leaves the most-recently-added one selected.
$q4.Activate()
